$wb = $excel.ActiveWorkbook

function Set-TextCell($range, [string]$text) {
    # Force literal text storage (avoids Excel's automatic date/number
    # inference for values like "2026-02-01"), then restore the cell's
    # style to the default "Normal" so no stray formatting is introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- Proximity sheet: append two new sensor events (rows 34-35) ---
$proximity = $wb.Worksheets.Item("Proximity")

Set-TextCell $proximity.Range("A34") "2026-02-01"
$proximity.Range("B34").Value = "15:12:46"
$proximity.Range("C34").Value = "15:00"
$proximity.Range("D34").Value = "Living Room Main Door"
$proximity.Range("E34").Value = "ENTER"
$proximity.Range("F34").Value = "User ENTERED Living Room Main Door"

Set-TextCell $proximity.Range("A35") "2026-02-01"
$proximity.Range("B35").Value = "15:12:47"
$proximity.Range("C35").Value = "15:00"
$proximity.Range("D35").Value = "Living Room Main Door"
$proximity.Range("E35").Value = "EXIT"
$proximity.Range("F35").Value = "User EXITED Living Room Main Door"

# --- Camera sheet: append one new image-captured event (row 7) ---
$camera = $wb.Worksheets.Item("Camera")

Set-TextCell $camera.Range("A7") "2026-02-01"
$camera.Range("B7").Value = "15:12:47"
$camera.Range("C7").Value = "15:00"
$camera.Range("D7").Value = "Living Room Main Door"
$camera.Range("E7").Value = "Image Captured"
$camera.Range("F7").Value = "Active"
